# Acta_nro_8.docx edits
# All edits are addressed via Tables/Cells (never via $d.Paragraphs — mixing
# document-level Paragraphs lookups with Table/Cell access corrupts the
# paragraph index cache in this host), and every replacement range is
# located by searching the owning cell's own Range.Text for the exact
# literal substring, then converting that in-cell offset into an absolute
# document Range. This keeps every edit minimal (touches only the
# characters that actually changed) and immune to any earlier edit shifting
# later offsets, since each cell is re-queried right before it is used.

$d = $word.ActiveDocument

function Set-SubstringInCell($cell, [string]$oldText, [string]$newText) {
    $cellStart = $cell.Range.Start
    $cellText = $cell.Range.Text
    $idx = $cellText.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Substring '$oldText' not found in cell text '$cellText'"
    }
    $targetStart = $cellStart + $idx
    $targetEnd = $targetStart + $oldText.Length
    $rng = $d.Range($targetStart, $targetEnd)
    $rng.Text = $newText
}

# --- Table 2: main meeting-info table (Fecha / Lugar / Hora Inicio / Hora Término) ---
$infoTable = $d.Tables.Item(2)

# "Lugar" value cell (row 2, col 3): "---" -> "Vía on-line"
Set-SubstringInCell $infoTable.Cell(2, 3) "---" "Vía on-line"

# "Hora Término" value cell (row 4, col 3): "--:-- hrs" -> "19:00 hrs"
Set-SubstringInCell $infoTable.Cell(4, 3) "--" "19"
Set-SubstringInCell $infoTable.Cell(4, 3) "--" "00"

# --- Table 6: "Próxima Reunión" footer table ---
$nextMeetingTable = $d.Tables.Item(6)

# Fecha: "9-10-2019" -> "11-10-2019"
Set-SubstringInCell $nextMeetingTable.Cell(2, 2) "9-10-2019" "11-10-2019"

# Objetivo de la Reunión bullet item: "---" -> new objective text
Set-SubstringInCell $nextMeetingTable.Cell(2, 3) "---" "Realizar puntos de  Formulación, Metodologías de trabajo y Plan de trabajo"

# Hora: "15:40 hrs" -> "17:00 hrs"
Set-SubstringInCell $nextMeetingTable.Cell(3, 2) "15" "17"
Set-SubstringInCell $nextMeetingTable.Cell(3, 2) "40" "00"

# Lugar: "U.B.B." -> "Vía on-line"
Set-SubstringInCell $nextMeetingTable.Cell(4, 2) "U.B.B." "Vía on-line"

Write-Output "Done."
